# Add a new "Brand SOS of Segment" exclusion rule row to the
# exclusion_rules sheet (mirrors the formatting of the other
# "location_type" / "Primary Shelf" rows already in the sheet),
# then tidy up sheet selections / active sheet to match the
# post-edit authoring session.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("exclusion_rules")
$ws2 = $wb.Worksheets.Item("store_policy_old")
$ws3 = $wb.Worksheets.Item("store_policy")

# --- exclusion_rules: append the new rule in row 22 ---
$ws1.Cells.Item(22, 1).Value = "Brand SOS of Segment"
$ws1.Cells.Item(22, 2).Value = "Include"
$ws1.Cells.Item(22, 3).Value = "location_type"
$ws1.Cells.Item(22, 4).Value = "Primary Shelf"

# Match the wrap-text styling used by column D on the other rows
$ws1.Cells.Item(22, 4).WrapText = $true

# --- selections left behind by the editing session ---
$ws2.Range("A6").Select()
$ws3.Range("A4").Select()

# exclusion_rules ends up the active sheet/tab with A22 selected
$ws1.Activate()
$ws1.Range("A22").Select()
